# "Updated our progress on the burndownchart"
#
# Sprint1 sheet: a few backlog items' remaining-work estimates (columns
# H:K, i.e. days 2-5 of the sprint) went down by one day each. The SUM
# row (row 28) and the chart-feeder cells (E33:E36), plus the burndown
# chart itself, all recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")
$ws.Activate()

# "Create a Car Class" (row 18): 3 -> 2
$ws.Range("H18:K18").Value2 = 2

# "MoreMaps" (row 20): 2 -> 1
$ws.Range("H20:K20").Value2 = 1

# "Class Info" (row 22): 2 -> 1
$ws.Range("H22:K22").Value2 = 1

# "Class Player" (row 23): 4 -> 3
$ws.Range("H23:K23").Value2 = 3

# The author also scrolled down a bit and left the selection on a
# different cell before saving.
$ws.Range("G26").Select() | Out-Null
